$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.227.61'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '1.895.93'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.17'
$ws.Range('E5').Value = '  +3.24%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5180'
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4006'
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08412'
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.69'
$ws.Range('E10').Value = '  +0.55%  '
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.37'
$ws.Range('E12').Value = '  +12.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.435'
$ws.Range('E13').Value = '  +1.89%  '
$ws.Range('D14').Value = '1.892.30'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.330'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.33'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06646'
$ws.Range('E19').Value = '  -1.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.21'
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.950'
$ws.Range('E22').Value = '  -1.50%  '
$ws.Range('D23').Value = '30.220.79'
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('E24').Value = '  +0.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.237'
$ws.Range('E25').Value = '  +1.24%  '
$ws.Range('D26').Value = '2.109.32'
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.58'
$ws.Range('E27').Value = '  +2.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.75'
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.352'
$ws.Range('E29').Value = '  -3.78%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '129.52'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.092'
$ws.Range('E31').Value = '  +2.34%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.090'
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.747'
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02498'
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06549'
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.259'
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.221'
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.78'
$ws.Range('E40').Value = '  +4.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.743'
$ws.Range('E41').Value = '  -3.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6499'
$ws.Range('E42').Value = '  -0.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.227'
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6098'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.19'
$ws.Range('E45').Value = '  +0.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.699'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.55'
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.01'
